$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range("Z1").Formula = '=""&"' + $text + '"'
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range("D2").Value = '27.562.55'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.646.83'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  -0.07%  '
Set-TextValue "D5" "212.56"
$ws.Range("E5").Value = '  -0.64%  '
Set-TextValue "D6" "0.537"
$ws.Range("E6").Value = '  +5.25%  '
$ws.Range("E7").Value = '  -0.04%  '
Set-TextValue "D8" "23.55"
$ws.Range("E8").Value = '  -2.04%  '
$ws.Range("E9").Value = '  -1.84%  '
$ws.Range("E10").Value = '  -1.30%  '
Set-TextValue "D11" "0.0888"
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '1.880.11'
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").Value = '1.657.69'
$ws.Range("E13").Value = '  -0.62%  '
Set-TextValue "D14" "0.584"
$ws.Range("E14").Value = '  +3.34%  '
Set-TextValue "D15" "4.03"
$ws.Range("E15").Value = '  -2.47%  '
Set-TextValue "D16" "64.45"
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("D17").Value = '27.530.27'
$ws.Range("E17").Value = '  -0.18%  '
Set-TextValue "D18" "230.89"
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("E19").Value = '  -0.65%  '
Set-TextValue "D20" "7.54"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("E22").Value = '  -3.34%  '
Set-TextValue "D23" "9.72"
$ws.Range("E23").Value = '  +3.68%  '
$ws.Range("E24").Value = '  -2.12%  '
Set-TextValue "D25" "149.04"
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("E26").Value = '  -2.69%  '
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("E29").Value = '  -4.21%  '
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("E31").Value = '  -3.24%  '
Set-TextValue "D32" "3.30"
$ws.Range("E32").Value = '  -0.77%  '
Set-TextValue "D33" "3.18"
$ws.Range("E33").Value = '  +2.26%  '
$ws.Range("D34").Value = '1.424.74'
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("E35").Value = '  +1.53%  '
$ws.Range("E36").Value = '  -0.03%  '
Set-TextValue "D37" "0.568"
$ws.Range("E37").Value = '  -0.58%  '
Set-TextValue "D38" "0.884"
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D42" "0.819"
$ws.Range("E42").Value = '  +3.17%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D43" "5.54"
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("E44").Value = '  +1.17%  '
Set-TextValue "D45" "65.06"
$ws.Range("E45").Value = '  -6.56%  '
$ws.Range("D46").Value = '1.789.32'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("E47").Value = '  -1.50%  '
Set-TextValue "D48" "88.12"
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("E50").Value = '  -2.91%  '
Set-TextValue "D51" "7.77"
$ws.Range("E51").Value = '  -0.59%  '

$ws.Range("Z1").ClearContents()

